$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: spreadsheet row number plus the new Coin/Link/Price/Volume(1h) values.
$updates = @(
    @{ Row = 2; D = '29.554.06'; E = '  +1.21%  ' },
    @{ Row = 3; D = '1.884.66'; E = '  +1.76%  ' },
    @{ Row = 4; D = '1.001'; E = '  -0.05%  ' },
    @{ Row = 5; D = '0.7164'; E = '  +2.85%  ' },
    @{ Row = 6; D = '242.82'; E = '  +2.22%  ' },
    @{ Row = 7; D = '1.001'; E = '  -0.03%  ' },
    @{ Row = 8; D = '0.07958'; E = '  +1.38%  ' },
    @{ Row = 9; D = '0.3123'; E = '  +3.69%  ' },
    @{ Row = 10; D = '25.41'; E = '  +7.84%  ' },
    @{ Row = 11; D = '0.08296'; E = '  +2.43%  ' },
    @{ Row = 12; D = '0.7326'; E = '  +4.09%  ' },
    @{ Row = 13; D = '1.887.15'; E = '  +2.09%  ' },
    @{ Row = 14; D = '5.297'; E = '  +2.25%  ' },
    @{ Row = 15; D = '91.62'; E = '  +2.45%  ' },
    @{ Row = 16; D = '29.549.50'; E = '  +1.06%  ' },
    @{ Row = 17; D = '5.961'; E = '  +2.71%  ' },
    @{ Row = 18; D = '247.57'; E = '  +5.18%  ' },
    @{ Row = 19; D = '0.000007904'; E = '  +1.38%  ' },
    @{ Row = 20; D = '13.41'; E = '  +1.80%  ' },
    @{ Row = 21; D = '2.138.28'; E = '  +1.66%  ' },
    @{ Row = 22; D = '1.000'; E = '  -0.10%  ' },
    @{ Row = 23; D = '7.998'; E = '  +6.46%  ' },
    @{ Row = 24; D = '1.001'; E = '  -0.12%  ' },
    @{ Row = 25; E = '  +14.64%  ' },
    @{ Row = 26; D = '163.55'; E = '  +0.56%  ' },
    @{ Row = 27; D = '9.090'; E = '  +2.77%  ' },
    @{ Row = 28; D = '18.40'; E = '  +2.17%  ' },
    @{ Row = 29; D = '1.357'; E = '  -3.67%  ' },
    @{ Row = 30; D = '1.500'; E = '  +1.95%  ' },
    @{ Row = 31; E = '  +2.07%  ' },
    @{ Row = 32; D = '4.129'; E = '  +2.98%  ' },
    @{ Row = 33; D = '0.05309'; E = '  +3.20%  ' },
    @{ Row = 34; D = '1.959'; E = '  +2.61%  ' },
    @{ Row = 35; E = '  +3.92%  ' },
    @{ Row = 36; D = '0.7306'; E = '  +2.87%  ' },
    @{ Row = 37; D = '2.680'; E = '  -0.06%  ' },
    @{ Row = 38; D = '0.01877'; E = '  +1.90%  ' },
    @{ Row = 39; D = '1.230.98'; E = '  +6.86%  ' },
    @{ Row = 40; D = '2.726'; E = '  +0.67%  ' },
    @{ Row = 41; D = '0.9142'; E = '  -1.12%  ' },
    @{ Row = 42; D = '74.76'; E = '  +6.99%  ' },
    @{ Row = 43; D = '6.231'; E = '  +4.35%  ' },
    @{ Row = 44; D = '1.001'; E = '  +0.03%  ' },
    @{ Row = 45; D = '102.39'; E = '  -0.33%  ' },
    @{ Row = 46; D = '2.035.81'; E = '  +2.14%  ' },
    @{ Row = 47; D = '0.5289'; E = '  -0.21%  ' },
    @{ Row = 48; E = '  +3.41%  ' },
    @{ Row = 49; D = '2.945'; E = '  +11.12%  ' },
    @{ Row = 50; B = 'BabyDogeCoin'; C = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D = '0.00000000120'; E = '  +2.89%  ' },
    @{ Row = 51; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '9.370'; E = '  +2.80%  ' }
)

# Price (D) and Volume(1h) (E) are written with a leading apostrophe so Excel
# keeps them as literal text -- this source column mixes values that look
# numeric ('1.000', '29.554.06') with ones that aren't valid numbers at all,
# and all of them must round-trip as the exact text shown in the source data.
foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = "'$($u.D)" }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = "'$($u.E)" }
}
